$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competicao") had been mistakenly truncated/saved as 55 for
# every data row; restore the correct competition id value of 255.
# xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 175
}

$ws.Range("B2:B$lastRow").Value = 255
